$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right above current row 256
# (pushes the former rows 256-274 down to 258-276).
$ws.Rows.Item(256).Insert()
$ws.Rows.Item(256).Insert()

# ---- New row 256: Venus / Primera ----
$ws.Range("A256").Value = 7
$ws.Range("B256").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C256").Value = "Ñuble"
$ws.Range("D256").Value = 44585
$ws.Range("E256").Value = 16
$ws.Range("F256").Value = "Fruta"
$ws.Range("G256").Value = 100103
$ws.Range("H256").Value = "Frutos de hueso (carozo)"
$ws.Range("I256").Value = 100103006
$ws.Range("J256").Value = "Nectarín"
$ws.Range("K256").Value = "Venus"
$ws.Range("L256").Value = "Primera"
$ws.Range("M256").Value = 120
$ws.Range("N256").Value = 12000
$ws.Range("O256").Value = 13000
$ws.Range("P256").Value = 12500
$ws.Range("Q256").Value = "$/caja 16 kilos empedrada"
$ws.Range("R256").Value = "Región de O'Higgins"
$ws.Range("S256").Value = 781
$ws.Range("T256").Value = 16

# ---- New row 257: Venus / Segunda ----
$ws.Range("A257").Value = 7
$ws.Range("B257").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C257").Value = "Ñuble"
$ws.Range("D257").Value = 44585
$ws.Range("E257").Value = 16
$ws.Range("F257").Value = "Fruta"
$ws.Range("G257").Value = 100103
$ws.Range("H257").Value = "Frutos de hueso (carozo)"
$ws.Range("I257").Value = 100103006
$ws.Range("J257").Value = "Nectarín"
$ws.Range("K257").Value = "Venus"
$ws.Range("L257").Value = "Segunda"
$ws.Range("M257").Value = 100
$ws.Range("N257").Value = 10000
$ws.Range("O257").Value = 11000
$ws.Range("P257").Value = 10500
$ws.Range("Q257").Value = "$/caja 16 kilos empedrada"
$ws.Range("R257").Value = "Región de O'Higgins"
$ws.Range("S257").Value = 656
$ws.Range("T257").Value = 16
